# Locate the shape that holds the "GitHub: <link>" text (slide 5, shape 3 in the
# original deck) by scanning for the text instead of hard-coding indices, so the
# script is resilient to any shape re-ordering.
$p = $ppt.ActivePresentation

$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "GitHub:*") {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find the 'GitHub: <link>' text box"
}

$tr = $targetShape.TextFrame.TextRange

# Old text:   "GitHub: https://github.com/pramodcgupta/MMC_Usecases/blob/main/EDA/Fire_Dept_Calls_Analytics.ipynb"
# New text:   "GitHub: https://nbviewer.org/github/pramodcgupta/MMC_Usecases/blob/main/Fire_Dept_Calls_Analytics.ipynb"
# ... with "GitHub" and ": " becoming separate runs, and a new blank paragraph
# inserted right after that line (pushing the existing blank/line-break
# paragraph further down).

$oldUrl = "https://github.com/pramodcgupta/MMC_Usecases/blob/main/EDA/Fire_Dept_Calls_Analytics.ipynb"
$newUrl = "https://nbviewer.org/github/pramodcgupta/MMC_Usecases/blob/main/Fire_Dept_Calls_Analytics.ipynb"

# 1) Insert a brand-new empty paragraph between paragraph 1 ("GitHub: <link>")
#    and the paragraph that follows it, by inserting a paragraph break right
#    before the second paragraph.
$secondParagraph = $tr.Paragraphs(2, 1)
$null = $secondParagraph.InsertBefore([char]13)

# 2) Split the leading bold run "GitHub: " into two runs: "GitHub" and ": ".
#    Nudging the formatting of just the ": " substring forces the engine to
#    materialize it as its own run while leaving "GitHub" untouched.
$colonSpace = $tr.Characters(7, 2)
$colonSpace.Font.Bold = $true

# 3) Swap the hyperlink run's visible text for the new nbviewer URL (the
#    run keeps its existing hyperlink relationship/formatting).
$linkText = $tr.Characters(9, $oldUrl.Length)
$linkText.Text = $newUrl
